# Add three new "DAILY ROUTINE" entries (rows 7-9) to Sheet1, mirroring the
# formatting of the existing rows (date column uses the same date style as
# the rows above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 12-Oct-2019 ---------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats - reuse the date style
$ws.Range("A7").Value = 43750
$ws.Range("B7").Value = "Register,login page were designed. The overall website flow was drawn"
$ws.Range("C7").Value = "Web pages related to the project was done"

# --- Row 8: 12-Nov-2019 (date only) ---------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 43781

# --- Row 9: 12-Dec-2019 ----------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 43811
$ws.Range("C9").Value = "Data flow diagram for the project was drawn "
$ws.Range("B9").Value = "Data flow,front screens for admin,faculty and student was drawn using html and css"

# Match the final selection recorded in the saved workbook
$ws.Range("B9").Select()
